$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.284.07"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.931.24"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7497"
$ws.Range("E5").Value = "  +4.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.34"
$ws.Range("E6").Value = "  -2.60%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "27.75"
$ws.Range("E8").Value = "  -0.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3179"
$ws.Range("E9").Value = "  -0.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07143"
$ws.Range("E10").Value = "  +0.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7802"
$ws.Range("E11").Value = "  -1.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08039"
$ws.Range("E12").Value = "  +0.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.918.39"
$ws.Range("E13").Value = "  -0.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.392"
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.97"
$ws.Range("E15").Value = "  -1.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.52"
$ws.Range("E16").Value = "  -1.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.294.66"
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.057"
$ws.Range("E18").Value = "  +5.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "251.48"
$ws.Range("E19").Value = "  -1.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007947"
$ws.Range("E20").Value = "  -1.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.172.85"
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.668"
$ws.Range("E24").Value = "  -2.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.530"
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.43"
$ws.Range("E26").Value = "  -0.53%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1296"
$ws.Range("E28").Value = "  +2.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.188"
$ws.Range("E29").Value = "  -3.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.369"
$ws.Range("E30").Value = "  +0.79%  "
$ws.Range("E31").Value = "  +1.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.407"
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.143"
$ws.Range("E33").Value = "  +0.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05230"
$ws.Range("E34").Value = "  +1.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.321"
$ws.Range("E35").Value = "  +4.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7547"
$ws.Range("E36").Value = "  +1.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.788"
$ws.Range("E37").Value = "  +1.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01953"
$ws.Range("E38").Value = "  -0.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.800"
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "78.45"
$ws.Range("E40").Value = "  +0.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.492"
$ws.Range("E41").Value = "  +1.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4518"
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.973"
$ws.Range("E43").Value = "  -0.84%  "
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8398"
$ws.Range("E45").Value = "  -0.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.04"
$ws.Range("E46").Value = "  +3.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.685"
$ws.Range("E47").Value = "  +3.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "101.69"
$ws.Range("E48").Value = "  +1.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.78"
$ws.Range("E49").Value = "  +3.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1210"
$ws.Range("E50").Value = "  +6.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "953.21"
$ws.Range("E51").Value = "  +3.57%  "
